$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ntng1"
$ws.Range("C2").Value = "Lrrc4c"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.025196
$ws.Range("H2").Value = 0.075588
$ws.Range("I2").Value = 0.01673902202018037
$ws.Range("J2").Value = 0.01673902202018037
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.03537
$ws.Range("N2").Value = 0.10611
$ws.Range("O2").Value = 0.01357924161260117
$ws.Range("P2").Value = 0.01357924161260117
$ws.Range("Q2").Value = 0.00089118252
$ws.Range("R2").Value = 0.00802064268
$ws.Range("S2").Value = 0.0002273032243706805
$ws.Range("T2").Value = 0.0002273032243706805

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ntng1"
$ws.Range("C3").Value = "Lrrc4c"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.025196
$ws.Range("H3").Value = 0.075588
$ws.Range("I3").Value = 0.01673902202018037
$ws.Range("J3").Value = 0.01673902202018037
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.011401
$ws.Range("N3").Value = 0.034203
$ws.Range("O3").Value = 0.004377069087511052
$ws.Range("P3").Value = 0.004377069087511052
$ws.Range("Q3").Value = 0.000287259596
$ws.Range("R3").Value = 0.002585336364
$ws.Range("S3").Value = 0.00007326785583969829
$ws.Range("T3").Value = 0.00007326785583969829

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ntng1"
$ws.Range("C4").Value = "Lrrc4c"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.025196
$ws.Range("H4").Value = 0.075588
$ws.Range("I4").Value = 0.01673902202018037
$ws.Range("J4").Value = 0.01673902202018037
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.55794
$ws.Range("N4").Value = 7.67382
$ws.Range("O4").Value = 0.9820436892998877
$ws.Range("P4").Value = 0.9820436892998878
$ws.Range("Q4").Value = 0.06444985623999999
$ws.Range("R4").Value = 0.58004870616
$ws.Range("S4").Value = 0.01643845093996999
$ws.Range("T4").Value = 0.01643845093996999

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ntng1"
$ws.Range("C5").Value = "Lrrc4c"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.03066833333333334
$ws.Range("H5").Value = 0.092005
$ws.Range("I5").Value = 0.02037457957568258
$ws.Range("J5").Value = 0.02037457957568258
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.03537
$ws.Range("N5").Value = 0.10611
$ws.Range("O5").Value = 0.01357924161260117
$ws.Range("P5").Value = 0.01357924161260117
$ws.Range("Q5").Value = 0.00108473895
$ws.Range("R5").Value = 0.00976265055
$ws.Range("S5").Value = 0.0002766713388133627
$ws.Range("T5").Value = 0.0002766713388133627

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ntng1"
$ws.Range("C6").Value = "Lrrc4c"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.03066833333333334
$ws.Range("H6").Value = 0.092005
$ws.Range("I6").Value = 0.02037457957568258
$ws.Range("J6").Value = 0.02037457957568258
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.011401
$ws.Range("N6").Value = 0.034203
$ws.Range("O6").Value = 0.004377069087511052
$ws.Range("P6").Value = 0.004377069087511052
$ws.Range("Q6").Value = 0.0003496496683333333
$ws.Range("R6").Value = 0.003146847015
$ws.Range("S6").Value = 0.00008918094243175428
$ws.Range("T6").Value = 0.00008918094243175427

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ntng1"
$ws.Range("C7").Value = "Lrrc4c"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.03066833333333334
$ws.Range("H7").Value = 0.092005
$ws.Range("I7").Value = 0.02037457957568258
$ws.Range("J7").Value = 0.02037457957568258
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.55794
$ws.Range("N7").Value = 7.67382
$ws.Range("O7").Value = 0.9820436892998877
$ws.Range("P7").Value = 0.9820436892998878
$ws.Range("Q7").Value = 0.07844775656666667
$ws.Range("R7").Value = 0.7060298091
$ws.Range("S7").Value = 0.02000872729443746
$ws.Range("T7").Value = 0.02000872729443746

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Ntng1"
$ws.Range("C8").Value = "Lrrc4c"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.449361
$ws.Range("H8").Value = 4.348083
$ws.Range("I8").Value = 0.962886398404137
$ws.Range("J8").Value = 0.962886398404137
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.03537
$ws.Range("N8").Value = 0.10611
$ws.Range("O8").Value = 0.01357924161260117
$ws.Range("P8").Value = 0.01357924161260117
$ws.Range("Q8").Value = 0.05126389856999999
$ws.Range("R8").Value = 0.46137508713
$ws.Range("S8").Value = 0.01307526704941712
$ws.Range("T8").Value = 0.01307526704941712

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Ntng1"
$ws.Range("C9").Value = "Lrrc4c"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.449361
$ws.Range("H9").Value = 4.348083
$ws.Range("I9").Value = 0.962886398404137
$ws.Range("J9").Value = 0.962886398404137
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.011401
$ws.Range("N9").Value = 0.034203
$ws.Range("O9").Value = 0.004377069087511052
$ws.Range("P9").Value = 0.004377069087511052
$ws.Range("Q9").Value = 0.016524164761
$ws.Range("R9").Value = 0.148717482849
$ws.Range("S9").Value = 0.0042146202892396
$ws.Range("T9").Value = 0.0042146202892396

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Ntng1"
$ws.Range("C10").Value = "Lrrc4c"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.449361
$ws.Range("H10").Value = 4.348083
$ws.Range("I10").Value = 0.962886398404137
$ws.Range("J10").Value = 0.962886398404137
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.55794
$ws.Range("N10").Value = 7.67382
$ws.Range("O10").Value = 0.9820436892998877
$ws.Range("P10").Value = 0.9820436892998878
$ws.Range("Q10").Value = 3.70737847634
$ws.Range("R10").Value = 33.36640628706
$ws.Range("S10").Value = 0.9455965110654803
$ws.Range("T10").Value = 0.9455965110654804
